# Initial version of holding conversion
#
# Insert a new "Preferred Conversion" column between "Option Type" and
# "Quantity" (i.e. before column J), which shifts the existing
# Quantity/Price/Grant Date columns one column to the right (J->K, K->L,
# L->M), and populate the new column's header plus the two known
# conversion-ratio data points (rows 3 and 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J — like selecting column J in the Excel UI and
# choosing Insert; carries formatting/styles/column-widths of the old J
# column onward with the shifted data (old J:L -> new K:M).
$ws.Columns.Item(10).Insert() | Out-Null

# Header for the newly inserted column.
$ws.Range("J1").Value = "Preferred Conversion"

# Known preferred-conversion ratios for the two "Preferred" option-type rows.
$ws.Range("J3").Value = 2
$ws.Range("J5").Value = 3

# Match the author's final cursor position/selection.
$ws.Range("J2").Select() | Out-Null
